$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking values keep their exact formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.034.20"
$ws.Range("E2").Value = "  -6.14%  "
$ws.Range("D3").Value = "2.429.01"
$ws.Range("E3").Value = "  -9.47%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "533.15"
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("D6").Value = "145.91"
$ws.Range("E6").Value = "  -7.49%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  -4.43%  "
$ws.Range("D9").Value = "0.0982"
$ws.Range("E9").Value = "  -7.15%  "
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("D11").Value = "5.37"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  -5.46%  "
$ws.Range("D13").Value = "2.860.92"
$ws.Range("E13").Value = "  -9.40%  "
$ws.Range("D14").Value = "23.75"
$ws.Range("E14").Value = "  -9.33%  "
$ws.Range("D15").Value = "58.985.54"
$ws.Range("E15").Value = "  -6.04%  "
$ws.Range("D16").Value = "0.0000137"
$ws.Range("E16").Value = "  -7.09%  "
$ws.Range("D17").Value = "2.478.71"
$ws.Range("E17").Value = "  -7.68%  "
$ws.Range("D18").Value = "11.06"
$ws.Range("E18").Value = "  -7.02%  "
$ws.Range("D19").Value = "4.34"
$ws.Range("E19").Value = "  -5.38%  "
$ws.Range("D20").Value = "322.83"
$ws.Range("E20").Value = "  -6.26%  "
$ws.Range("D21").Value = "0.966"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").Value = "5.67"
$ws.Range("E22").Value = "  -9.92%  "
$ws.Range("D23").Value = "0.462"
$ws.Range("E23").Value = "  -8.40%  "
$ws.Range("D24").Value = "59.88"
$ws.Range("E24").Value = "  -5.27%  "
$ws.Range("D25").Value = "0.159"
$ws.Range("E25").Value = "  -5.29%  "
$ws.Range("D26").Value = "0.973"
$ws.Range("E26").Value = "  -2.43%  "
$ws.Range("D27").Value = "7.64"
$ws.Range("E27").Value = "  -6.32%  "
$ws.Range("D28").Value = "1.27"
$ws.Range("E28").Value = "  -5.57%  "
$ws.Range("D29").Value = "6.73"
$ws.Range("E29").Value = "  -4.98%  "
$ws.Range("D30").Value = "1.81"
$ws.Range("E30").Value = "  -6.17%  "
$ws.Range("D31").Value = "0.0₃0755"
$ws.Range("E31").Value = "  -11.77%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "155.41"
$ws.Range("E33").Value = "  -6.31%  "
$ws.Range("D34").Value = "4.47"
$ws.Range("E34").Value = "  -7.62%  "
$ws.Range("D35").Value = "18.24"
$ws.Range("E35").Value = "  -6.63%  "
$ws.Range("D36").Value = "1.34"
$ws.Range("E36").Value = "  -7.23%  "
$ws.Range("D37").Value = "1.72"
$ws.Range("E37").Value = "  -3.02%  "
$ws.Range("D38").Value = "310.24"
$ws.Range("E38").Value = "  -8.40%  "
$ws.Range("D39").Value = "5.70"
$ws.Range("E39").Value = "  -8.20%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "36.69"
$ws.Range("E40").Value = "  -4.07%  "
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").Value = "0.837"
$ws.Range("E41").Value = "  -11.15%  "
$ws.Range("D42").Value = "3.68"
$ws.Range("E42").Value = "  -7.08%  "
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "10.73"
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.0934"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.575"
$ws.Range("E46").Value = "  -6.86%  "
$ws.Range("D47").Value = "0.0522"
$ws.Range("E47").Value = "  -7.62%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0227"
$ws.Range("E48").Value = "  -5.51%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "18.34"
$ws.Range("E49").Value = "  -9.76%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.962.53"
$ws.Range("E50").Value = "  -6.22%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "18.64"
$ws.Range("E51").Value = "  -10.47%  "
